# Updated cryptos list on Sun Nov 17 09:54:48 UTC 2024 with GitHub Actions
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in this
# sheet. Where the new Price value would otherwise be auto-parsed as a
# number by Excel (e.g. "237.32"), it is entered with a leading apostrophe
# to force text entry, then the cell style is reset to "Normal" so the
# quote-prefix formatting left behind by that trick doesn't change the
# cell's style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.177.70"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "3.135.62"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'237.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.68%  "

$ws.Range("D6").Value = "'634.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("E7").Value = "  +4.67%  "

$ws.Range("D8").Value = "'0.366"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "3.135.92"
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("D11").Value = "'0.723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "

$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "'36.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.25%  "

$ws.Range("D14").Value = "'0.0000249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "'5.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").Value = "91.012.47"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").Value = "3.721.02"
$ws.Range("E17").Value = "  +0.78%  "

$ws.Range("D18").Value = "3.153.04"
$ws.Range("E18").Value = "  +1.83%  "

$ws.Range("E19").Value = "  -3.88%  "

$ws.Range("D20").Value = "'0.0000214"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").Value = "'14.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").Value = "'446.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("D23").Value = "'5.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.51%  "

$ws.Range("D24").Value = "'9.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("D25").Value = "'5.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.79%  "

$ws.Range("D26").Value = "'90.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").Value = "'12.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "'9.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.97%  "

$ws.Range("E31").Value = "  -3.70%  "

$ws.Range("D32").Value = "'0.997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.10%  "

$ws.Range("D33").Value = "'0.200"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.48%  "

$ws.Range("D34").Value = "'26.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.76%  "

$ws.Range("D35").Value = "'3.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.53%  "

$ws.Range("D36").Value = "'514.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.18%  "

$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").Value = "'7.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("E39").Value = "  +4.09%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +5.37%  "

$ws.Range("D42").Value = "'22.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "

$ws.Range("D43").Value = "'0.0852"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'3.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +48.38%  "

$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").Value = "'151.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").Value = "'0.697"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.70%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'45.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.27%  "

$ws.Range("E51").Value = "  +2.93%  "
